$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.135
$ws.Range("A3").Value = -21.445
$ws.Range("C5").Value = -13.058
$ws.Range("D5").Value = -8.309000000000001
$ws.Range("E7").Value = 13.078
$ws.Range("D9").Value = -7.752000000000001
$ws.Range("D11").Value = -8.284000000000001
$ws.Range("E11").Value = 12.87
$ws.Range("A14").Value = -20.964
$ws.Range("E19").Value = 12.753
$ws.Range("A21").Value = -20.828
$ws.Range("D21").Value = -7.936
$ws.Range("E21").Value = 13.293
$ws.Range("A23").Value = -21.701
$ws.Range("A25").Value = -22.077
